# Progress Report template updated
#
# The document ends with a table (Decision List) followed by three empty
# "Body"-styled paragraphs right before the final section break (sectPr):
#   1) <w:p><w:pPr><w:pStyle w:val="Body"/><w:widowControl w:val="0"/></w:pPr></w:p>
#   2) <w:p><w:pPr><w:pStyle w:val="Body"/></w:pPr></w:p>
#   3) <w:p><w:pPr><w:pStyle w:val="Body"/></w:pPr><w:r/></w:p>
#
# The edit removes the second and third (trailing) empty paragraphs, keeping
# only the first one (the paragraph with widowControl explicitly set to 0)
# immediately before the sectPr.

$d = $word.ActiveDocument

# Locate the end of the last table in the document (the "Decision List"
# table). The paragraph immediately following it is the one we must keep;
# the two paragraphs after that one must be removed.
$tableCount = $d.Tables.Count
$lastTable = $d.Tables.Item($tableCount)
$keepParaEnd = $lastTable.Range.End + 1

# Deleting the paragraph mark right after the paragraph we keep merges the
# (empty) following paragraph away. Doing this twice removes both trailing
# empty paragraphs, leaving only the first one intact.
$d.Range($keepParaEnd, $keepParaEnd + 1).Delete()
$d.Range($keepParaEnd, $keepParaEnd + 1).Delete()
